$wb = $excel.ActiveWorkbook

$alc = $wb.Worksheets.Item("ALC")
$arm = $wb.Worksheets.Item("ARM")
$bsm = $wb.Worksheets.Item("BSM")
$crp = $wb.Worksheets.Item("CRP")
$cul = $wb.Worksheets.Item("CUL")
$wvr = $wb.Worksheets.Item("WVR")

# ALC!row58
$alc.Range("H58").Value = 1953.6364
$alc.Range("I58").Value = 236
$alc.Range("J58").Value = 3385
$alc.Range("K58").Value = 708
$alc.Range("L58").Value = 10155
$alc.Range("M58").Value = -558
$alc.Range("N58").Value = -10455

# ALC!row62
$alc.Range("H62").Value = 20641.842
$alc.Range("I62").Value = 5282.8335
$alc.Range("J62").Value = 46971.57
$alc.Range("K62").Value = 5282.8335
$alc.Range("L62").Value = 46971.57
$alc.Range("M62").Value = -4658.8335
$alc.Range("N62").Value = -48219.57

# ALC!row65
$alc.Range("H65").Value = 20641.842
$alc.Range("I65").Value = 5282.8335
$alc.Range("J65").Value = 46971.57
$alc.Range("K65").Value = 26414.1675
$alc.Range("L65").Value = 234857.85
$alc.Range("M65").Value = -23294.1675
$alc.Range("N65").Value = -241097.85

# ALC!row116
$alc.Range("H116").Value = 4637.1875
$alc.Range("I116").Value = 2745
$alc.Range("J116").Value = 8800
$alc.Range("K116").Value = 2745
$alc.Range("L116").Value = 8800
$alc.Range("M116").Value = 697
$alc.Range("N116").Value = -15684

# ALC!row132
$alc.Range("H132").Value = 5352.806
$alc.Range("I132").Value = 4787.8696
$alc.Range("J132").Value = 6590.2856
$alc.Range("K132").Value = 14363.6088
$alc.Range("L132").Value = 19770.8568
$alc.Range("M132").Value = -11833.6088
$alc.Range("N132").Value = -24830.8568

# ALC!row137
$alc.Range("H137").Value = 1600.7407
$alc.Range("I137").Value = 1658.697
$alc.Range("J137").Value = 1509.6666
$alc.Range("K137").Value = 4976.090999999999
$alc.Range("L137").Value = 4528.9998
$alc.Range("M137").Value = -2426.090999999999
$alc.Range("N137").Value = -9628.9998

# ALC!row138
$alc.Range("H138").Value = 1718.2037
$alc.Range("I138").Value = 797.5909
$alc.Range("J138").Value = 5768.9
$alc.Range("K138").Value = 2392.7727
$alc.Range("L138").Value = 17306.7
$alc.Range("M138").Value = 2747.2273
$alc.Range("N138").Value = -27586.7

# ARM!row2
$arm.Range("H2").Value = 2626.15
$arm.Range("I2").Value = 1711.7368
$arm.Range("J2").Value = 20000
$arm.Range("K2").Value = 1711.7368
$arm.Range("L2").Value = 20000
$arm.Range("M2").Value = -1598.7368
$arm.Range("N2").Value = -20226

# ARM!row102
$arm.Range("H102").Value = 3000
$arm.Range("I102").Value = 3000
$arm.Range("J102").Value = 0
$arm.Range("K102").Value = 3000
$arm.Range("L102").Value = 0
$arm.Range("M102").Value = -1378
$arm.Range("N102").ClearContents()

# ARM!row116
$arm.Range("H116").Value = 2626.15
$arm.Range("I116").Value = 1711.7368
$arm.Range("J116").Value = 20000
$arm.Range("K116").Value = 1711.7368
$arm.Range("L116").Value = 20000
$arm.Range("M116").Value = 582.2632000000001
$arm.Range("N116").Value = -24588

# ARM!row122
$arm.Range("H122").Value = 1677.6
$arm.Range("I122").Value = 1568.2667
$arm.Range("J122").Value = 1841.6
$arm.Range("K122").Value = 4704.800099999999
$arm.Range("L122").Value = 5524.799999999999
$arm.Range("M122").Value = -2254.800099999999
$arm.Range("N122").Value = -10424.8

# BSM!row3
$bsm.Range("H3").Value = 2626.15
$bsm.Range("I3").Value = 1711.7368
$bsm.Range("J3").Value = 20000
$bsm.Range("K3").Value = 1711.7368
$bsm.Range("L3").Value = 20000
$bsm.Range("M3").Value = -1597.7368
$bsm.Range("N3").Value = -20228

# BSM!row53
$bsm.Range("H53").Value = 37890
$bsm.Range("J53").Value = 37890
$bsm.Range("L53").Value = 37890
$bsm.Range("N53").Value = -39038

# BSM!row64
$bsm.Range("H64").Value = 2673.6667
$bsm.Range("I64").Value = 1208.4
$bsm.Range("J64").Value = 10000
$bsm.Range("K64").Value = 1208.4
$bsm.Range("L64").Value = 10000
$bsm.Range("M64").Value = -983.4000000000001
$bsm.Range("N64").Value = -10450

# BSM!row67
$bsm.Range("H67").Value = 2673.6667
$bsm.Range("I67").Value = 1208.4
$bsm.Range("J67").Value = 10000
$bsm.Range("K67").Value = 1208.4
$bsm.Range("L67").Value = 10000
$bsm.Range("M67").Value = -428.4000000000001
$bsm.Range("N67").Value = -11560

# CRP!row16
$crp.Range("H16").Value = 2676.9092
$crp.Range("I16").Value = 3007.5454
$crp.Range("J16").Value = 2346.2727
$crp.Range("K16").Value = 3007.5454
$crp.Range("L16").Value = 2346.2727
$crp.Range("M16").Value = -2720.5454
$crp.Range("N16").Value = -2920.2727

# CRP!row31
$crp.Range("H31").Value = 4445978
$crp.Range("I31").Value = 1134.2407
$crp.Range("J31").Value = 15875577
$crp.Range("K31").Value = 1134.2407
$crp.Range("L31").Value = 15875577
$crp.Range("M31").Value = -839.2407000000001
$crp.Range("N31").Value = -15876167

# CRP!row34
$crp.Range("H34").Value = 4445978
$crp.Range("I34").Value = 1134.2407
$crp.Range("J34").Value = 15875577
$crp.Range("K34").Value = 1134.2407
$crp.Range("L34").Value = 15875577
$crp.Range("M34").Value = -932.2407000000001
$crp.Range("N34").Value = -15875981

# CRP!row41
$crp.Range("H41").Value = 15314.143
$crp.Range("I41").Value = 5000
$crp.Range("J41").Value = 17033.166
$crp.Range("K41").Value = 5000
$crp.Range("L41").Value = 17033.166
$crp.Range("M41").Value = -4572
$crp.Range("N41").Value = -17889.166

# CRP!row50
$crp.Range("H50").Value = 18740
$crp.Range("I50").Value = 0
$crp.Range("J50").Value = 18740
$crp.Range("K50").Value = 0
$crp.Range("L50").Value = 18740
$crp.Range("M50").ClearContents()
$crp.Range("N50").Value = -19990

# CRP!row51
$crp.Range("H51").Value = 19475
$crp.Range("I51").Value = 18500
$crp.Range("J51").Value = 19800
$crp.Range("K51").Value = 18500
$crp.Range("L51").Value = 19800
$crp.Range("M51").Value = -17764
$crp.Range("N51").Value = -21272

# CRP!row59
$crp.Range("H59").Value = 46864.707
$crp.Range("J59").Value = 47981.816
$crp.Range("L59").Value = 47981.816
$crp.Range("N59").Value = -50271.816

# CRP!row60
$crp.Range("H60").Value = 11110.833
$crp.Range("I60").Value = 5000
$crp.Range("J60").Value = 11376.521
$crp.Range("K60").Value = 5000
$crp.Range("L60").Value = 11376.521
$crp.Range("M60").Value = -4489
$crp.Range("N60").Value = -12398.521

# CRP!row61
$crp.Range("H61").Value = 19475
$crp.Range("I61").Value = 18500
$crp.Range("J61").Value = 19800
$crp.Range("K61").Value = 18500
$crp.Range("L61").Value = 19800
$crp.Range("M61").Value = -18152
$crp.Range("N61").Value = -20496

# CRP!row68
$crp.Range("H68").Value = 29800
$crp.Range("J68").Value = 29800
$crp.Range("L68").Value = 29800
$crp.Range("N68").Value = -31298

# CRP!row71
$crp.Range("H71").Value = 29800
$crp.Range("J71").Value = 29800
$crp.Range("L71").Value = 89400
$crp.Range("N71").Value = -96888

# CRP!row74
$crp.Range("H74").Value = 33950
$crp.Range("J74").Value = 33950
$crp.Range("L74").Value = 33950
$crp.Range("N74").Value = -35698

# CRP!row77
$crp.Range("H77").Value = 33950
$crp.Range("J77").Value = 33950
$crp.Range("L77").Value = 101850
$crp.Range("N77").Value = -110586

# CRP!row113
$crp.Range("H113").Value = 2676.9092
$crp.Range("I113").Value = 3007.5454
$crp.Range("J113").Value = 2346.2727
$crp.Range("K113").Value = 3007.5454
$crp.Range("L113").Value = 2346.2727
$crp.Range("M113").Value = -837.5454
$crp.Range("N113").Value = -6686.2727

# CUL!row117
$cul.Range("H117").Value = 28571784
$cul.Range("I117").Value = 447
$cul.Range("J117").Value = 142857140
$cul.Range("K117").Value = 1341
$cul.Range("L117").Value = 428571420
$cul.Range("M117").Value = 2101
$cul.Range("N117").Value = -428578304

# CUL!row131
$cul.Range("H131").Value = 323411.44
$cul.Range("I131").Value = 667150.4
$cul.Range("J131").Value = 1156.1875
$cul.Range("K131").Value = 2001451.2
$cul.Range("L131").Value = 3468.5625
$cul.Range("M131").Value = -1996411.2
$cul.Range("N131").Value = -13548.5625

# WVR!row122
$wvr.Range("H122").Value = 4276.9585
$wvr.Range("I122").Value = 3402.238
$wvr.Range("K122").Value = 10206.714
$wvr.Range("M122").Value = -7756.714
